# Apply "Correct Forecast output" restructuring:
#  - Sheet1 renamed "Sales vs PO", gains a new "Order Week" column (C) holding
#    the original ds dates, the ds column itself shifted one week later, and
#    PO_Requested_Qty (now column D) zeroed out (actuals moved to sheet 2).
#  - New sheet "Weekly Growth" holding the PO qty + week-over-week growth %.
#  - New sheet "Volume Insights" holding summary stats.
#  - New sheet "Prediction Info" holding the next-week prediction.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Sales vs PO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Original column-A (ds) values, before we overwrite them - needed for the
# new "Order Week" column C.
$origDs = @(45565, 45572, 45579, 45586, 45593, 45600, 45607, 45614, 45621, 45628, 45635, 45642, 45649)

# 1) Give the new header cell D1 the same style as the existing header C1
#    (bold / centered / bordered), then copy column-A's date style down into
#    the new column C (format only, so no new style entries are created).
$ws1.Range("A1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)

$ws1.Range("A2:A14").Copy()
$ws1.Range("C2:C14").PasteSpecial(-4122)

# 2) Header row text
$ws1.Cells.Item(1, 3).Value = "Order Week"
$ws1.Cells.Item(1, 4).Value = "PO_Requested_Qty"

# 3) Fill column C ("Order Week") with the original ds values, and set the
#    new PO_Requested_Qty column (D) to 0 for every data row.
for ($i = 0; $i -lt $origDs.Length; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 3).Value = $origDs[$i]
    $ws1.Cells.Item($r, 4).Value = 0
}

# 4) Shift column A (ds) forward by one week (+6 -> net +7 vs the value that
#    used to sit there, landing on the following Wednesday in this data).
for ($i = 0; $i -lt $origDs.Length; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $origDs[$i] + 6
}

# ---------------------------------------------------------------------------
# Sheet 2: "Weekly Growth"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws2.Range("B1").PasteSpecial(-4122)
$ws2.Range("C1").PasteSpecial(-4122)

$ws1.Range("A2:A4").Copy()
$ws2.Range("A2:A4").PasteSpecial(-4122)

$ws2.Cells.Item(1, 1).Value = "ds"
$ws2.Cells.Item(1, 2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1, 3).Value = "Growth%"

$ws2.Cells.Item(2, 1).Value = 45572
$ws2.Cells.Item(2, 2).Value = 16
$ws2.Cells.Item(2, 3).Value = 0

$ws2.Cells.Item(3, 1).Value = 45586
$ws2.Cells.Item(3, 2).Value = 224
$ws2.Cells.Item(3, 3).Value = 1300

$ws2.Cells.Item(4, 1).Value = 45607
$ws2.Cells.Item(4, 2).Value = 128
$ws2.Cells.Item(4, 3).Value = -42.85714285714286

# ---------------------------------------------------------------------------
# Sheet 3: "Volume Insights"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws1.Range("A1").Copy()
$ws3.Range("A1").PasteSpecial(-4122)
$ws3.Range("B1").PasteSpecial(-4122)
$ws3.Range("C1").PasteSpecial(-4122)
$ws3.Range("D1").PasteSpecial(-4122)

$ws3.Cells.Item(1, 1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1, 2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1, 3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1, 4).Value = "Min_PO_Quantity"

$ws3.Cells.Item(2, 1).Value = 368
$ws3.Cells.Item(2, 2).Value = 122.6666666666667
$ws3.Cells.Item(2, 3).Value = 224
$ws3.Cells.Item(2, 4).Value = 16

# ---------------------------------------------------------------------------
# Sheet 4: "Prediction Info"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

$ws4.Cells.Item(1, 1).Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Cells.Item(2, 1).Value = 234.6666666666666

# Re-select A1 on the first sheet to match the original view state.
$ws1.Range("A1").Select() | Out-Null
